# Updated symbol list on Wed Feb 15 22:17:57 UTC 2023 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values to the
# cryptocurrency table on Sheet1. Source values are stored as plain text
# (e.g. "313.16", "5.62%") rather than numbers/percentages, so each cell is
# forced to Text format before the write and the format is cleared again
# afterwards to avoid leaving a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $range = $ws.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.ClearFormats()
}


Set-TextValue "D2" "313.16"
Set-TextValue "E2" "5.62%"
Set-TextValue "D3" "44.62"
Set-TextValue "E3" "7.11%"
Set-TextValue "D4" "5.154"
Set-TextValue "E4" "2.34%"
Set-TextValue "D5" "0.08040"
Set-TextValue "E5" "6.42%"
Set-TextValue "D6" "4.518"
Set-TextValue "E6" "2.84%"
Set-TextValue "D7" "1.657"
Set-TextValue "E7" "3.36%"
Set-TextValue "D8" "1.084"
Set-TextValue "E8" "16.69%"
Set-TextValue "D9" "0.1303"
Set-TextValue "E9" "8.98%"
Set-TextValue "D10" "0.1921"
Set-TextValue "E10" "4.19%"
Set-TextValue "D11" "0.09382"
Set-TextValue "E11" "4.43%"
Set-TextValue "D12" "0.04218"
Set-TextValue "E12" "2.46%"
Set-TextValue "D13" "0.1040"
Set-TextValue "E13" "-1.01%"
Set-TextValue "D14" "0.001314"
Set-TextValue "E14" "2.56%"
Set-TextValue "D15" "0.005895"
Set-TextValue "E15" "-1.64%"
Set-TextValue "D17" "3.387"
Set-TextValue "E17" "1.04%"
Set-TextValue "D18" "2.401"
Set-TextValue "E18" "-0.28%"
Set-TextValue "D19" "0.3375"
Set-TextValue "E19" "1.66%"
Set-TextValue "D20" "7.992"
Set-TextValue "E20" "1.11%"
Set-TextValue "D21" "0.1370"
Set-TextValue "E21" "-3.44%"
Set-TextValue "E22" "4.56%"
Set-TextValue "D23" "0.04209"
Set-TextValue "E23" "3.69%"
Set-TextValue "D24" "0.001273"
Set-TextValue "D25" "0.004577"
Set-TextValue "E25" "15.59%"
Set-TextValue "E26" "8.80%"
Set-TextValue "D38" "0.02672"
Set-TextValue "E38" "10.87%"
Set-TextValue "D39" "0.05421"
Set-TextValue "E39" "3.96%"
Set-TextValue "D40" "0.005626"
Set-TextValue "E40" "-13.20%"
Set-TextValue "D41" "0.007754"
Set-TextValue "E41" "-0.22%"
Set-TextValue "D42" "0.1417"
Set-TextValue "E42" "6.41%"
Set-TextValue "D43" "0.007342"
Set-TextValue "E43" "-2.86%"
Set-TextValue "D44" "0.007925"
Set-TextValue "E44" "1.07%"
Set-TextValue "D45" "0.3126"
Set-TextValue "E45" "-2.59%"
Set-TextValue "D46" "0.00006797"
Set-TextValue "E46" "0.24%"
Set-TextValue "D47" "0.00000000745"
Set-TextValue "E47" "-0.69%"
Set-TextValue "D48" "0.05902"
Set-TextValue "E48" "27.55%"
Set-TextValue "D49" "0.003973"
Set-TextValue "E49" "-5.45%"
Set-TextValue "D50" "0.00002086"
Set-TextValue "E50" "-0.69%"
Set-TextValue "D51" "0.0001987"
Set-TextValue "E51" "-0.69%"
